$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 35.995988
$ws.Range("H2").Value = 107.987964
$ws.Range("I2").Value = 0.5613901502831141
$ws.Range("J2").Value = 0.561390150283114
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.214110666666667
$ws.Range("N2").Value = 21.642332
$ws.Range("O2").Value = 0.4688823795981188
$ws.Range("P2").Value = 0.4688823795981188
$ws.Range("Q2").Value = 259.6790409880053
$ws.Range("R2").Value = 2337.111368892048
$ws.Range("S2").Value = 0.2632259495476921
$ws.Range("T2").Value = 0.2632259495476921

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 35.995988
$ws.Range("H3").Value = 107.987964
$ws.Range("I3").Value = 0.5613901502831141
$ws.Range("J3").Value = 0.561390150283114
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.110350666666666
$ws.Range("N3").Value = 21.331052
$ws.Range("O3").Value = 0.4621384803214003
$ws.Range("P3").Value = 0.4621384803214003
$ws.Range("Q3").Value = 255.9440972731254
$ws.Range("R3").Value = 2303.496875458128
$ws.Range("S3").Value = 0.2594399909192409
$ws.Range("T3").Value = 0.2594399909192409

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 35.995988
$ws.Range("H4").Value = 107.987964
$ws.Range("I4").Value = 0.5613901502831141
$ws.Range("J4").Value = 0.561390150283114
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.061296333333333
$ws.Range("N4").Value = 3.183889
$ws.Range("O4").Value = 0.06897914008048092
$ws.Range("P4").Value = 0.06897914008048092
$ws.Range("Q4").Value = 38.20241007911066
$ws.Range("R4").Value = 343.821690711996
$ws.Range("S4").Value = 0.03872420981618117
$ws.Range("T4").Value = 0.03872420981618116

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 20.666474
$ws.Range("H5").Value = 61.999422
$ws.Range("I5").Value = 0.3223124461726698
$ws.Range("J5").Value = 0.3223124461726698
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.214110666666667
$ws.Range("N5").Value = 21.642332
$ws.Range("O5").Value = 0.4688823795981188
$ws.Range("P5").Value = 0.4688823795981188
$ws.Range("Q5").Value = 149.0902305257893
$ws.Range("R5").Value = 1341.812074732104
$ws.Range("S5").Value = 0.151126626735532
$ws.Range("T5").Value = 0.151126626735532

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 20.666474
$ws.Range("H6").Value = 61.999422
$ws.Range("I6").Value = 0.3223124461726698
$ws.Range("J6").Value = 0.3223124461726698
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.110350666666666
$ws.Range("N6").Value = 21.331052
$ws.Range("O6").Value = 0.4621384803214003
$ws.Range("P6").Value = 0.4621384803214003
$ws.Range("Q6").Value = 146.9458771835493
$ws.Range("R6").Value = 1322.512894651944
$ws.Range("S6").Value = 0.1489529840629108
$ws.Range("T6").Value = 0.1489529840629108

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 20.666474
$ws.Range("H7").Value = 61.999422
$ws.Range("I7").Value = 0.3223124461726698
$ws.Range("J7").Value = 0.3223124461726698
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.061296333333333
$ws.Range("N7").Value = 3.183889
$ws.Range("O7").Value = 0.06897914008048092
$ws.Range("P7").Value = 0.06897914008048092
$ws.Range("Q7").Value = 21.93325307912866
$ws.Range("R7").Value = 197.399277712158
$ws.Range("S7").Value = 0.02223283537422706
$ws.Range("T7").Value = 0.02223283537422706

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.456917333333333
$ws.Range("H8").Value = 22.370752
$ws.Range("I8").Value = 0.116297403544216
$ws.Range("J8").Value = 0.116297403544216
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.214110666666667
$ws.Range("N8").Value = 21.642332
$ws.Range("O8").Value = 0.4688823795981188
$ws.Range("P8").Value = 0.4688823795981188
$ws.Range("Q8").Value = 53.79502687485155
$ws.Range("R8").Value = 484.155241873664
$ws.Range("S8").Value = 0.05452980331489471
$ws.Range("T8").Value = 0.05452980331489471

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.456917333333333
$ws.Range("H9").Value = 22.370752
$ws.Range("I9").Value = 0.116297403544216
$ws.Range("J9").Value = 0.116297403544216
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.110350666666666
$ws.Range("N9").Value = 21.331052
$ws.Range("O9").Value = 0.4621384803214003
$ws.Range("P9").Value = 0.4621384803214003
$ws.Range("Q9").Value = 53.02129713234488
$ws.Range("R9").Value = 477.191674191104
$ws.Range("S9").Value = 0.05374550533924863
$ws.Range("T9").Value = 0.05374550533924863

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.456917333333333
$ws.Range("H10").Value = 22.370752
$ws.Range("I10").Value = 0.116297403544216
$ws.Range("J10").Value = 0.116297403544216
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.061296333333333
$ws.Range("N10").Value = 3.183889
$ws.Range("O10").Value = 0.06897914008048092
$ws.Range("P10").Value = 0.06897914008048092
$ws.Range("Q10").Value = 7.913999023836443
$ws.Range("R10").Value = 71.225991214528
$ws.Range("S10").Value = 0.008022094890072697
$ws.Range("T10").Value = 0.008022094890072697
